$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.842.57"
$ws.Range("E2").Value = "  +3.94%  "
$ws.Range("D3").Value = "'2.280.89"
$ws.Range("E3").Value = "  +4.76%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "'252.58"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "'0.639"
$ws.Range("E6").Value = "  +4.32%  "
$ws.Range("D7").Value = "'72.62"
$ws.Range("E7").Value = "  +9.25%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.653"
$ws.Range("E9").Value = "  +13.26%  "
$ws.Range("D10").Value = "'38.66"
$ws.Range("E10").Value = "  +6.18%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "'59.85"
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").Value = "'0.0969"
$ws.Range("E12").Value = "  +3.53%  "
$ws.Range("D13").Value = "'7.38"
$ws.Range("E13").Value = "  +7.47%  "
$ws.Range("D14").Value = "'0.106"
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("D15").Value = "'2.622.76"
$ws.Range("E15").Value = "  +4.99%  "
$ws.Range("D16").Value = "'15.00"
$ws.Range("E16").Value = "  +4.64%  "
$ws.Range("D17").Value = "'0.886"
$ws.Range("E17").Value = "  +4.48%  "
$ws.Range("D18").Value = "'2.278.83"
$ws.Range("E18").Value = "  +4.67%  "
$ws.Range("D19").Value = "'42.811.93"
$ws.Range("E19").Value = "  +4.24%  "
$ws.Range("E20").Value = "  +7.00%  "
$ws.Range("D21").Value = "'6.34"
$ws.Range("E21").Value = "  +4.37%  "
$ws.Range("D22").Value = "'73.48"
$ws.Range("E22").Value = "  +2.52%  "
$ws.Range("D23").Value = "'236.83"
$ws.Range("E23").Value = "  +2.82%  "
$ws.Range("E24").Value = "  +3.55%  "
$ws.Range("D25").Value = "'3.88"
$ws.Range("E25").Value = "  +2.24%  "
$ws.Range("D26").Value = "'11.63"
$ws.Range("E26").Value = "  +0.85%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("D29").Value = "'3.68"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("E30").Value = "  +4.56%  "
$ws.Range("D31").Value = "'168.09"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "'21.07"
$ws.Range("E32").Value = "  +4.02%  "
$ws.Range("D33").Value = "'6.39"
$ws.Range("E33").Value = "  +9.16%  "
$ws.Range("D34").Value = "'0.130"
$ws.Range("E34").Value = "  +7.14%  "
$ws.Range("D35").Value = "'0.0805"
$ws.Range("E35").Value = "  +5.84%  "
$ws.Range("D36").Value = "'30.65"
$ws.Range("E36").Value = "  +24.20%  "
$ws.Range("E37").Value = "  +4.63%  "
$ws.Range("D38").Value = "'4.70"
$ws.Range("E38").Value = "  +18.61%  "
$ws.Range("D39").Value = "'4.78"
$ws.Range("E39").Value = "  +5.44%  "
$ws.Range("D40").Value = "'0.0312"
$ws.Range("E40").Value = "  +2.25%  "
$ws.Range("E41").Value = "  +4.98%  "
$ws.Range("D42").Value = "'13.28"
$ws.Range("E42").Value = "  +16.07%  "
$ws.Range("D43").Value = "'5.96"
$ws.Range("E43").Value = "  +8.53%  "
$ws.Range("D44").Value = "'0.213"
$ws.Range("E44").Value = "  +12.36%  "
$ws.Range("D45").Value = "'9.19"
$ws.Range("E45").Value = "  +7.98%  "
$ws.Range("E46").Value = "  -5.74%  "
$ws.Range("D47").Value = "'61.51"
$ws.Range("E47").Value = "  +1.05%  "
$ws.Range("E48").Value = "  +2.44%  "
$ws.Range("E49").Value = "  +3.97%  "
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("E51").Value = "  +4.93%  "
